$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "parisk"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "DIS"
$ws.Range("E3").Value = "OTH"
$ws.Range("F3").Value = "c8048836-24fe-4e27-95aa-c7cfb58ac155"
$ws.Range("G3").Value = "rkc_hGb0Z_annotated.xlsx"
$ws.Range("H3").Value = "The structure of the global policies used in the experiments should be mentioned somewhere."
